$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Pearson logo in the primary footer: image1.png -> image2.png
$ftrPrimary = $sec.Footers.Item(1)
$shpA = $ftrPrimary.Range.InlineShapes.Item(1)
$shapeA = $shpA.ConvertToShape()
$shapeA.Name = "image2.png"
$shapeA.ConvertToInlineShape() | Out-Null

# Pearson logo in the first-page footer: image1.png -> image2.png
$ftrFirst = $sec.Footers.Item(2)
$shpB = $ftrFirst.Range.InlineShapes.Item(1)
$shapeB = $shpB.ConvertToShape()
$shapeB.Name = "image2.png"
$shapeB.ConvertToInlineShape() | Out-Null

# BTec logo in the first-page header: image2.jpg -> image1.jpg
$hdrFirst = $sec.Headers.Item(2)
$shpC = $hdrFirst.Range.InlineShapes.Item(1)
$shapeC = $shpC.ConvertToShape()
$shapeC.Name = "image1.jpg"
$shapeC.ConvertToInlineShape() | Out-Null

Write-Output "done"
